$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Change 1: table column widths were nudged (total width preserved) ---
# old dxa: 440 / 2645 / 1276 / 4359  ->  new dxa: 440 / 2641 / 1287 / 4352
# dxa -> points is /20
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 2).Width = 132.05
    $t.Cell($r, 3).Width = 64.35
    $t.Cell($r, 4).Width = 217.6
}

# --- Change 2: the "No se si se hace por empresa" comment (row with "1" /
# "se hace por empresa") was split across three runs with proofErr spell
# markers around "se"; collapse it back into a single plain run. A
# Find/Replace across the whole phrase merges the runs and drops the
# now-unneeded proofErr markers. ---
$d.Content.Find.Execute("No se si se hace por empresa", $true, $false, $false, `
    $false, $false, $true, 1, $false, "No se si se hace por empresa", 2) | Out-Null

# Locate the "Ficha de un proveedor" row dynamically (status / comment cells).
$targetRow = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    if ($t.Cell($r, 2).Range.Text -like "*Ficha de un proveedor*") {
        $targetRow = $r
    }
}

# --- Change 3: status cell "No" -> "Imcompleta" (flagged as a spelling
# error by Word, hence the surrounding proofErr markers). ---
$cellStatus = $t.Cell($targetRow, 3)
$rngStatus = $cellStatus.Range
$rngStatus.End = $rngStatus.End - 1
$rngStatus.Text = ""
$xmlStatus = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:cnfStyle w:val="000000000000"/><w:rPr><w:lang w:val="es-VE"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t>Imcompleta</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngStatus.InsertXML($xmlStatus) | Out-Null

# --- Change 4: the (empty) comment cell gets a new numbered-list
# paragraph explaining what is missing. ---
$cellComment = $t.Cell($targetRow, 4)
$rngComment = $cellComment.Range
$rngComment.End = $rngComment.End - 1
$rngComment.Text = ""
$xmlComment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:cnfStyle w:val="000000000000"/><w:rPr><w:lang w:val="es-VE"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t xml:space="preserve"> No </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t>se</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t xml:space="preserve"> si es necesario agregar pa&#237;ses  a los que distribuye, formas de pago y de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t>envio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-VE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngComment.InsertXML($xmlComment) | Out-Null
